$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Column C (rows 2-11): 46078 -> 46079 ("Förändrad" date column)
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 3).Value = 46079
}

# 2. Rows 7-11: content in columns A (Beteckning), B (Datum), G (Area (ha))
#    is permuted between rows. Capture the "before" values first so the
#    subsequent writes don't clobber values we still need to read.
$rows = 7, 8, 9, 10, 11
$colA = @{}
$colB = @{}
$colG = @{}
foreach ($r in $rows) {
    $colA[$r] = $ws.Cells.Item($r, 1).Value()
    $colB[$r] = $ws.Cells.Item($r, 2).Value()
    $colG[$r] = $ws.Cells.Item($r, 7).Value()
}

# Permutation: new row <- old row
#   7  <- 11
#   8  <- 10
#   9  <- 7
#   10 <- 8
#   11 <- 9
$mapping = @{ 7 = 11; 8 = 10; 9 = 7; 10 = 8; 11 = 9 }

foreach ($newRow in $rows) {
    $oldRow = $mapping[$newRow]
    $ws.Cells.Item($newRow, 1).Value = $colA[$oldRow]
    $ws.Cells.Item($newRow, 2).Value = $colB[$oldRow]
    $ws.Cells.Item($newRow, 7).Value = $colG[$oldRow]
}
